$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 18
$ws.Range("H18").Value = 1360.5
$ws.Range("I18").Value = 1432.2
$ws.Range("J18").Value = 1002
$ws.Range("K18").Value = 1432.2
$ws.Range("L18").Value = 1002
$ws.Range("M18").Value = -1148.2
$ws.Range("N18").Value = -1570
# Row 40
$ws.Range("H40").Value = 3554.077
$ws.Range("J40").Value = 3109.182
$ws.Range("L40").Value = 3109.182
$ws.Range("N40").Value = -3459.182
# Row 54
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
# Row 62
$ws.Range("H62").Value = 49229.223
$ws.Range("I62").Value = 51599.5
$ws.Range("K62").Value = 51599.5
$ws.Range("M62").Value = -50975.5
# Row 65
$ws.Range("H65").Value = 49229.223
$ws.Range("I65").Value = 51599.5
$ws.Range("K65").Value = 257997.5
$ws.Range("M65").Value = -254877.5
# Row 100
$ws.Range("H100").Value = 2279.6
# Row 135
$ws.Range("H135").Value = 43479572
$ws.Range("I135").Value = 43479572
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 391316148
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -391313613
$ws.Range("N135").ClearContents()
# Row 137
$ws.Range("H137").Value = 2657.0667
$ws.Range("I137").Value = 2450.5386
$ws.Range("K137").Value = 7351.6158
$ws.Range("M137").Value = -4801.6158

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 2849.2666
$ws.Range("I5").Value = 293.25
$ws.Range("J5").Value = 5770.4287
$ws.Range("K5").Value = 293.25
$ws.Range("L5").Value = 5770.4287
$ws.Range("M5").Value = -181.25
$ws.Range("N5").Value = -5994.4287
# Row 32
$ws.Range("H32").Value = 3151.7856
$ws.Range("I32").Value = 3125.1282
$ws.Range("J32").Value = 3498.3333
$ws.Range("K32").Value = 3125.1282
$ws.Range("L32").Value = 3498.3333
$ws.Range("M32").Value = -2838.1282
$ws.Range("N32").Value = -4072.3333
# Row 45
$ws.Range("H45").Value = 1844.5834
$ws.Range("I45").Value = 1439.75
$ws.Range("J45").Value = 2654.25
$ws.Range("K45").Value = 1439.75
$ws.Range("L45").Value = 2654.25
$ws.Range("M45").Value = -1062.75
$ws.Range("N45").Value = -3408.25
# Row 88
$ws.Range("H88").Value = 5953670
# Row 91
$ws.Range("H91").Value = 5953670
# Row 132
$ws.Range("H132").Value = 58825670
$ws.Range("I132").Value = 66668836
$ws.Range("K132").Value = 200006508
$ws.Range("M132").Value = -200003978
# Row 133
$ws.Range("H133").Value = 68472.5
$ws.Range("J133").Value = 68472.5
$ws.Range("L133").Value = 68472.5
$ws.Range("N133").Value = -73532.5

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 2849.2666
$ws.Range("I4").Value = 293.25
$ws.Range("J4").Value = 5770.4287
$ws.Range("K4").Value = 293.25
$ws.Range("L4").Value = 5770.4287
$ws.Range("M4").Value = -178.25
$ws.Range("N4").Value = -6000.4287
# Row 20
$ws.Range("H20").Value = 10988.318
$ws.Range("I20").Value = 14828.786
$ws.Range("J20").Value = 4267.5
$ws.Range("K20").Value = 14828.786
$ws.Range("L20").Value = 4267.5
$ws.Range("M20").Value = -14581.786
$ws.Range("N20").Value = -4761.5
# Row 26
$ws.Range("H26").Value = 15471
$ws.Range("I26").Value = 15471
$ws.Range("K26").Value = 15471
$ws.Range("M26").Value = -15179
# Row 96
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("M96").ClearContents()
# Row 134
$ws.Range("H134").Value = 1655.7333
$ws.Range("I134").Value = 1655.7333
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4967.199900000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2432.199900000001
$ws.Range("N134").ClearContents()

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2644.8635
$ws.Range("I31").Value = 1979.4
$ws.Range("J31").Value = 4070.8572
$ws.Range("K31").Value = 1979.4
$ws.Range("L31").Value = 4070.8572
$ws.Range("M31").Value = -1684.4
$ws.Range("N31").Value = -4660.8572
# Row 34
$ws.Range("H34").Value = 2644.8635
$ws.Range("I34").Value = 1979.4
$ws.Range("J34").Value = 4070.8572
$ws.Range("K34").Value = 1979.4
$ws.Range("L34").Value = 4070.8572
$ws.Range("M34").Value = -1777.4
$ws.Range("N34").Value = -4474.8572
# Row 50
$ws.Range("H50").Value = 63991.168
$ws.Range("I50").Value = 30000
$ws.Range("J50").Value = 70789.39999999999
$ws.Range("K50").Value = 30000
$ws.Range("L50").Value = 70789.39999999999
$ws.Range("M50").Value = -29375
$ws.Range("N50").Value = -72039.39999999999
# Row 62
$ws.Range("H62").Value = 66673530
$ws.Range("I62").Value = 5988
$ws.Range("J62").Value = 166674830
$ws.Range("K62").Value = 5988
$ws.Range("L62").Value = 166674830
$ws.Range("M62").Value = -5364
$ws.Range("N62").Value = -166676078
# Row 65
$ws.Range("H65").Value = 66673530
$ws.Range("I65").Value = 5988
$ws.Range("J65").Value = 166674830
$ws.Range("K65").Value = 29940
$ws.Range("L65").Value = 833374150
$ws.Range("M65").Value = -26820
$ws.Range("N65").Value = -833380390
# Row 69
$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14251
# Row 72
$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -41256
# Row 105
$ws.Range("H105").Value = 1354.909
$ws.Range("I105").Value = 1379.3
$ws.Range("J105").Value = 1111
$ws.Range("K105").Value = 1379.3
$ws.Range("L105").Value = 1111
$ws.Range("M105").Value = 367.7
$ws.Range("N105").Value = -4605
# Row 135
$ws.Range("H135").Value = 52986.668
$ws.Range("J135").Value = 52986.668
$ws.Range("L135").Value = 52986.668
$ws.Range("N135").Value = -63126.668

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1388.8889
$ws.Range("I68").Value = 1243
$ws.Range("J68").Value = 1899.5
$ws.Range("K68").Value = 3729
$ws.Range("L68").Value = 5698.5
$ws.Range("M68").Value = -2918
$ws.Range("N68").Value = -7320.5
# Row 71
$ws.Range("H71").Value = 1388.8889
$ws.Range("I71").Value = 1243
$ws.Range("J71").Value = 1899.5
$ws.Range("K71").Value = 11187
$ws.Range("L71").Value = 17095.5
$ws.Range("M71").Value = -7131
$ws.Range("N71").Value = -25207.5
# Row 122
$ws.Range("H122").Value = 455.83334
$ws.Range("J122").Value = 345.33334
$ws.Range("L122").Value = 3108.00006
$ws.Range("N122").Value = -8008.00006
# Row 131
$ws.Range("H131").Value = 6399.6
$ws.Range("J131").Value = 6641.684
$ws.Range("L131").Value = 19925.052
$ws.Range("N131").Value = -30005.052
# Row 139
$ws.Range("H139").Value = 4777812
$ws.Range("I139").Value = 5573280.5
$ws.Range("K139").Value = 16719841.5
$ws.Range("M139").Value = -16714701.5
# Row 141
$ws.Range("H141").Value = 2256.5715
$ws.Range("I141").Value = 2256.5715
$ws.Range("K141").Value = 6769.7145
$ws.Range("M141").Value = -1589.7145

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2373.4707
$ws.Range("I80").Value = 2326.2
$ws.Range("K80").Value = 2326.2
$ws.Range("M80").Value = -1328.2
# Row 83
$ws.Range("H83").Value = 2373.4707
$ws.Range("I83").Value = 2326.2
$ws.Range("K83").Value = 11631
$ws.Range("M83").Value = -6639
# Row 98
$ws.Range("H98").Value = 10550
$ws.Range("J98").Value = 10550
$ws.Range("L98").Value = 10550
$ws.Range("N98").Value = -16540
# Row 102
$ws.Range("H102").Value = 1779.3448
$ws.Range("I102").Value = 894.6
$ws.Range("K102").Value = 894.6
$ws.Range("M102").Value = 727.4
# Row 122
$ws.Range("H122").Value = 31253740
$ws.Range("I122").Value = 2483.3333
$ws.Range("K122").Value = 7449.999899999999
$ws.Range("M122").Value = -4999.999899999999
# Row 132
$ws.Range("H132").Value = 3362.8108
$ws.Range("I132").Value = 2931
$ws.Range("J132").Value = 5213.4287
$ws.Range("K132").Value = 8793
$ws.Range("L132").Value = 15640.2861
$ws.Range("M132").Value = -6263
$ws.Range("N132").Value = -20700.2861

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2932961
$ws.Range("I22").Value = 353.33334
$ws.Range("J22").Value = 4032689
$ws.Range("K22").Value = 353.33334
$ws.Range("L22").Value = 4032689
$ws.Range("M22").Value = -58.33334000000002
$ws.Range("N22").Value = -4033279
# Row 27
$ws.Range("H27").Value = 2932961
$ws.Range("I27").Value = 353.33334
$ws.Range("J27").Value = 4032689
$ws.Range("K27").Value = 353.33334
$ws.Range("L27").Value = 4032689
$ws.Range("M27").Value = -246.33334
$ws.Range("N27").Value = -4032903
# Row 55
$ws.Range("H55").Value = 335.54544
$ws.Range("I55").Value = 339.1
$ws.Range("K55").Value = 339.1
$ws.Range("M55").Value = -166.1
# Row 96
$ws.Range("H96").Value = 200001
$ws.Range("J96").Value = 200001
$ws.Range("L96").Value = 200001
$ws.Range("N96").Value = -205493
# Row 136
$ws.Range("H136").Value = 2447
$ws.Range("I136").Value = 2184.8572
$ws.Range("K136").Value = 6554.571599999999
$ws.Range("M136").Value = -4004.571599999999

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 5724.5
$ws.Range("J96").Value = 7250.75
$ws.Range("L96").Value = 7250.75
$ws.Range("N96").Value = -9996.75
# Row 126
$ws.Range("H126").Value = 1824.1177
$ws.Range("I126").Value = 1615.7142
$ws.Range("K126").Value = 4847.142599999999
$ws.Range("M126").Value = -2377.142599999999
# Row 138
$ws.Range("H138").Value = 98000
$ws.Range("J138").Value = 98000
$ws.Range("L138").Value = 98000
$ws.Range("N138").Value = -108280
